$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the sprint date headers (row 4, E:K) to the new sprint dates ---
$ws.Range("E4").Value = 45215
$ws.Range("F4").Value = 45216
$ws.Range("G4").Value = 45217
$ws.Range("H4").Value = 45218
$ws.Range("I4").Value = 45219
$ws.Range("J4").Value = 45220
$ws.Range("K4").Value = 45221

# --- Remove the old task rows 6-14 (tasks 7-14, i.e. rows 11-19) ---
# Tasks 1-5 (rows 6-10) are kept/updated below; the remaining 9 task rows
# (old tasks 6-14, rows 11-19) are deleted, shifting the summary rows
# (Completed Effort / Remaining Effort / Ideal Burndown) from rows 20-22
# up to rows 11-13.
$ws.Rows("11:19").Delete()

# Clear the stray formatted-but-empty L11 cell (was L20) left outside the
# new used range.
$ws.Range("L11").Clear()

# --- Update the 5 remaining task rows with the new work breakdown ---
# Task 1
$ws.Range("C6").Value = "Set up the gitHub repository"
$ws.Range("D6").Value = 4
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 2

# Task 2
$ws.Range("C7").Value = "Set up the game"
$ws.Range("D7").Value = 3
$ws.Range("I7").Value = 3

# Task 3
$ws.Range("C8").Value = "Test the game"
$ws.Range("D8").Value = 5
$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 2

# Task 4
$ws.Range("C9").Value = "Prepare the structure of future scrums, burndown charts and sprints"
$ws.Range("D9").Value = 3
$ws.Range("I9").Value = 2

# Task 5
$ws.Range("C10").Value = "First quick look at the code"
$ws.Range("D10").Value = 5
$ws.Range("J10").Value = 2
